$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $CellRef, $Text)
    $c = $Sheet.Range($CellRef)
    $c.Value = "'" + $Text
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "261.27"
Set-TextValue $ws "E2" "1.72%"
Set-TextValue $ws "D3" "27.37"
Set-TextValue $ws "E3" "1.81%"
Set-TextValue $ws "D4" "4.708"
Set-TextValue $ws "E4" "1.45%"
Set-TextValue $ws "D5" "0.06079"
Set-TextValue $ws "E5" "3.13%"
Set-TextValue $ws "D6" "6.672"
Set-TextValue $ws "E6" "0.98%"
Set-TextValue $ws "D7" "0.8459"
Set-TextValue $ws "D8" "0.9250"
Set-TextValue $ws "E8" "0.35%"
Set-TextValue $ws "D9" "0.1403"
Set-TextValue $ws "E9" "1.95%"
Set-TextValue $ws "D10" "0.04889"
Set-TextValue $ws "E10" "14.05%"
Set-TextValue $ws "D11" "0.07112"
Set-TextValue $ws "E11" "1.59%"
Set-TextValue $ws "D12" "0.03097"
Set-TextValue $ws "E12" "1.26%"
Set-TextValue $ws "D13" "0.09070"
Set-TextValue $ws "D14" "0.001531"
Set-TextValue $ws "E14" "0.14%"
Set-TextValue $ws "D15" "0.0006069"
Set-TextValue $ws "E15" "-0.06%"
Set-TextValue $ws "D16" "0.006126"
Set-TextValue $ws "E16" "2.15%"
Set-TextValue $ws "D17" "3.451"
Set-TextValue $ws "E17" "-0.57%"
Set-TextValue $ws "D18" "3.141"
Set-TextValue $ws "E18" "-0.82%"
Set-TextValue $ws "D20" "0.3109"
Set-TextValue $ws "E20" "1.76%"
Set-TextValue $ws "D21" "0.1289"
Set-TextValue $ws "E21" "-0.53%"
Set-TextValue $ws "D22" "4.084"
Set-TextValue $ws "E22" "4.42%"
Set-TextValue $ws "D23" "0.04232"
Set-TextValue $ws "E23" "-0.57%"
Set-TextValue $ws "D24" "0.001221"
Set-TextValue $ws "E24" "0.16%"
Set-TextValue $ws "E25" "-8.82%"
Set-TextValue $ws "D26" "0.0001200"
Set-TextValue $ws "E26" "0.02%"
Set-TextValue $ws "D27" "0.0001576"
Set-TextValue $ws "E27" "3.38%"
Set-TextValue $ws "D40" "0.03868"
Set-TextValue $ws "E40" "2.39%"
Set-TextValue $ws "E41" "1.28%"
Set-TextValue $ws "D42" "0.004104"
Set-TextValue $ws "E42" "-34.12%"
Set-TextValue $ws "D43" "0.01635"
Set-TextValue $ws "E43" "15.77%"
Set-TextValue $ws "E44" "-8.39%"
Set-TextValue $ws "D45" "0.00005147"
Set-TextValue $ws "E45" "-4.16%"
Set-TextValue $ws "E46" "0.02%"
Set-TextValue $ws "D47" "0.05447"
Set-TextValue $ws "E47" "19.67%"
Set-TextValue $ws "E49" "0.02%"
Set-TextValue $ws "E50" "0.02%"

Write-Host "Applied all cell updates"
